# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet named "2022-Q1" right after "2021-Q4" (i.e. right
#    before the "总计" summary sheet) and populate it with that quarter's
#    fund holdings detail (same shape as the other quarterly sheets).
# 2. Insert a new top data row into the "总计" summary sheet for "2022-Q1"
#    (pushing the existing quarters down by one row) and fill in its values.
#
# Note: worksheet object variables in this host resolve by tab *position*,
# not by a stable handle, so every sheet reference below is re-fetched by
# name immediately before it is used (never cached across a structural
# change such as Worksheets.Add/Move/rename).

$wb = $excel.ActiveWorkbook

# --- Make room for the new sheet in the right tab position -----------------
# Rename the existing "总计" sheet out of the way, then add a fresh sheet
# right after it and name that one "总计". This leaves the *new* blank sheet
# in the "总计" slot, and the old (renamed) sheet sitting where "2022-Q1"
# belongs, directly after "2021-Q4".
$wb.Worksheets.Item("总计").Name = "2022-Q1"
$placeholder = $wb.Worksheets.Add($null, $wb.Worksheets.Item("2022-Q1"))
$placeholder.Name = "总计"

# --- Build the "2022-Q1" sheet ----------------------------------------------
# This sheet currently still holds the old "总计" table's data (it was only
# renamed above) - wipe it before writing the quarter's fund holdings.
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Cells.Clear()

# Columns B & C (fund code / fund name) and D-G (size/position figures kept
# as formatted text, same as every other quarterly sheet) must stay text -
# mark them "@" (Text) before writing so numeric-looking strings (fund code
# "008099", ratios like "61.82") keep their literal form instead of being
# parsed into numbers.
$q1.Range("B1:G2").NumberFormat = "@"

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"
$q1.Range("B1:H1").Font.Bold = $true
$q1.Range("B1:H1").HorizontalAlignment = -4108
$q1.Range("B1:H1").VerticalAlignment = -4160
$q1.Range("B1:H1").Borders.LineStyle = 1

$q1.Cells.Item(2, 1).Value = 0
$q1.Range("A2").Font.Bold = $true
$q1.Range("A2").HorizontalAlignment = -4108
$q1.Range("A2").VerticalAlignment = -4160
$q1.Range("A2").Borders.LineStyle = 1
$q1.Cells.Item(2, 2).Value = "008099"
$q1.Cells.Item(2, 3).Value = "广发价值领先混合"
$q1.Cells.Item(2, 4).Value = "61.82"
$q1.Cells.Item(2, 5).Value = "83.88"
$q1.Cells.Item(2, 6).Value = "5.24"
$q1.Cells.Item(2, 7).Value = "3.2394"
$q1.Cells.Item(2, 8).Value = 3

# --- Insert the new "2022-Q1" row into the "总计" table --------------------
$total = $wb.Worksheets.Item("总计")

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"
$total.Range("B1:D1").Font.Bold = $true
$total.Range("B1:D1").HorizontalAlignment = -4108
$total.Range("B1:D1").VerticalAlignment = -4160
$total.Range("B1:D1").Borders.LineStyle = 1

$rows = @(
    @(0, "2022-Q1", 1, 3.24),
    @(1, "2021-Q4", 3, 0.59),
    @(2, "2021-Q3", 3, 4.04),
    @(3, "2021-Q2", 1, 4.1),
    @(4, "2021-Q1", 1, 1.97),
    @(5, "2020-Q4", 1, 0.24)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Range($total.Cells.Item($r, 1), $total.Cells.Item($r, 1)).Font.Bold = $true
    $total.Range($total.Cells.Item($r, 1), $total.Cells.Item($r, 1)).HorizontalAlignment = -4108
    $total.Range($total.Cells.Item($r, 1), $total.Cells.Item($r, 1)).VerticalAlignment = -4160
    $total.Range($total.Cells.Item($r, 1), $total.Cells.Item($r, 1)).Borders.LineStyle = 1
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}

Write-Output "done"
